$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> becomes the "15612 / WER" subject record
$ws.Range("A2").Value = 15612
$ws.Range("B2").Value = "WER"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = 89
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = 78
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = 167
$ws.Range("L2").Value = 20.875

# Row 3 -> becomes the "15611 / SAD" subject record
$ws.Range("A3").Value = 15611
$ws.Range("B3").Value = "SAD"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = 54
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = 56
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = 110
$ws.Range("L3").Value = 13.75
